$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix two small typos (rows 61 and 62) ---
$ws.Range("C61").Value = "A symptom that is related to a depression diagnosis."
$ws.Range("C62").Value = "Symptom severity that is associated with a depression symptom."

# --- Split the old row 96 ("post-traumatic stress disorder severity") into
#     a "symptom" row (96) and a new "symptom severity" row (97) ---

# 1) Insert a brand new row at position 97 (shifts old rows 97-120 down to 98-121)
$ws.Rows.Item(97).Insert()

# 2) Update row 96 in place: it now describes the post-traumatic stress SYMPTOM
$ws.Range("A96").Value = "GMHO:0000174"
$ws.Range("B96").Value = "post-traumatic stress symptom"
$ws.Range("C96").Value = "An anxiety symptom that related to a post-traumatic stress disorder diagnosis."
$ws.Range("D96").Value = "anxiety symptom"
$ws.Range("G96").Value = ""

# 3) Populate the newly inserted row 97: the post-traumatic stress SYMPTOM SEVERITY
$ws.Range("A97").Value = "GMHO:0000173"
$ws.Range("B97").Value = "post-traumatic stress symptom severity"
$ws.Range("C97").Value = "An anxiety symptom severity relating to a post-traumatic stress symptom."
$ws.Range("D97").Value = "anxiety disorder severity"
$ws.Range("G97").Value = "PTSD severity"
$ws.Range("S97").Value = "Proposed"
$ws.Range("V97").Value = "PS"
